$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.295.58"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.568.46"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "'207.80"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").Value = "'1.00"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "'0.478"
$ws.Range("E7").Value = "  -4.72%  "
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'17.88"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "1.785.46"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "1.567.31"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "'0.507"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "25.295.18"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "'59.49"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("E18").Value = "  -2.54%  "
$ws.Range("D20").Value = "'186.24"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").Value = "'4.15"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'9.30"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").Value = "'0.130"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'139.84"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("E27").Value = "  -6.97%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.46"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'14.87"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").Value = "'1.16"
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("D31").Value = "'0.0465"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").Value = "'3.05"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").Value = "'3.01"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("D36").Value = "1.088.83"
$ws.Range("E36").Value = "  -3.33%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("D40").Value = "'0.497"
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = "  -8.00%  "
$ws.Range("D42").Value = "'0.762"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'93.43"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").Value = "'5.07"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D45").Value = "1.698.34"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").Value = "0.0₆0105"
$ws.Range("E46").Value = "  -7.98%  "
$ws.Range("D47").Value = "'52.84"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("E49").Value = "  -3.21%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("E51").Value = "  -0.70%  "
